# Commit: Mon, Jun 15, 2020  1:05:03 AM
#
# The three "Component three" summary tables (on slides 14, 15 and 16)
# had their table style switched from the deck's default custom style
# ({41DAA5FE-7D41-4150-B446-7ED891656FCD}) to PowerPoint's built-in
# "No Style, No Grid" table style ({B8B55C0F-EDF5-4B7B-AAF2-388ACD2A277C}).

$p = $ppt.ActivePresentation

$targetStyleId = "{B8B55C0F-EDF5-4B7B-AAF2-388ACD2A277C}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
